$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.031.16"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.555.34"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.84"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3918"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3220"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.49"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07336"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.09"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.685"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.767"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.556.64"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001120"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06624"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.13"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.383"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.91"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.42"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.040.94"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.347"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.533"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.61"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.83"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.863"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.733.04"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.81"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.080"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.781"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.679"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -13.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.371"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08228"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06247"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02288"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.199"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2106"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.216"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.19%  "

$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.74"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5914"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.40"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.722"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5704"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.926"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.23"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.148"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06874"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.88%  "
